$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1546, 248.0625, 1297.9375)
    3  = @(1574, 262.3541666666667, 1311.645833333333)
    4  = @(1552, 252.1458333333333, 1299.854166666667)
    5  = @(1664, 359.75, 1304.25)
    6  = @(1496, 295.1041666666667, 1200.895833333333)
    7  = @(1520, 354.6875, 1165.3125)
    8  = @(1384, 243.2291666666667, 1140.770833333333)
    9  = @(1382, 255.4375, 1126.5625)
    10 = @(1406, 251.2291666666667, 1154.770833333333)
    11 = @(1420, 270.1041666666667, 1149.895833333333)
    12 = @(1530, 255.4375, 1274.5625)
    13 = @(1480, 242.8333333333333, 1237.166666666667)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}
